$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $Value) {
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '63.854.04'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.134.77'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue 'D5' '598.82'
$ws.Range('E5').Value = '  -2.52%  '
Set-TextValue 'D6' '139.64'
$ws.Range('E6').Value = '  -3.77%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.128.63'
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('E10').Value = '  -2.66%  '
Set-TextValue 'D11' '5.34'
$ws.Range('E11').Value = '  -1.41%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('E13').Value = '  -1.77%  '
Set-TextValue 'D14' '34.47'
$ws.Range('E14').Value = '  -3.01%  '
$ws.Range('D15').Value = '3.649.57'
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('D17').Value = '63.819.62'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '3.126.34'
$ws.Range('E18').Value = '  -1.06%  '
Set-TextValue 'D19' '6.76'
Set-TextValue 'D20' '482.19'
$ws.Range('E20').Value = '  +0.57%  '
Set-TextValue 'D21' '14.50'
$ws.Range('E21').Value = '  -0.52%  '
Set-TextValue 'D22' '0.704'
$ws.Range('E22').Value = '  -2.30%  '
Set-TextValue 'D23' '7.67'
$ws.Range('E23').Value = '  -3.54%  '
Set-TextValue 'D24' '87.66'
$ws.Range('E24').Value = '  +4.54%  '
Set-TextValue 'D25' '13.05'
$ws.Range('E25').Value = '  -5.28%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -2.45%  '
Set-TextValue 'D28' '8.12'
$ws.Range('E28').Value = '  -5.97%  '
Set-TextValue 'D29' '6.91'
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('E31').Value = '  +2.62%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D32' '1.00'
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.111'
$ws.Range('E33').Value = '  -7.41%  '
$ws.Range('E34').Value = '  -3.18%  '
$ws.Range('E35').Value = '  -1.99%  '
$ws.Range('E36').Value = '  +0.17%  '
Set-TextValue 'D37' '52.56'
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('E38').Value = '  -6.22%  '
Set-TextValue 'D39' '0.0394'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D40' '2.85'
$ws.Range('E40').Value = '  -10.50%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D41' '427.23'
$ws.Range('E41').Value = '  -7.33%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '2.875.00'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('E45').Value = '  -3.00%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D46' '2.38'
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D47' '2.15'
$ws.Range('E47').Value = '  -6.79%  '
Set-TextValue 'D48' '0.999'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  +0.18%  '
Set-TextValue 'D50' '25.54'
$ws.Range('E50').Value = '  -3.66%  '
Set-TextValue 'D51' '120.43'
$ws.Range('E51').Value = '  +0.67%  '
